$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Update the first three summary rows to "0M" ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Insert 10 new rows after row 3 (before the old row 4), carrying the
#     per-iteration values that used to be crammed into a single
#     tab-separated row further down the table ---
$refRow = $t.Rows.Item(4)
$newValues = @("138", "0.00003", "0.00012", "0.00003", "0.00001", "0.00004", "0.00004", "0.00004", "0.00531", "100.0")
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# --- Collapse the three tab-separated summary rows (now shifted down by
#     the 10 inserted rows) down to a single value each ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.01"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "147"
